$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22787.503246
$ws.Range("D2").Value = 1494.907896416667
$ws.Range("E2").Value = 4417.314503233333

$ws.Range("B3").Value = 21200.1327983
$ws.Range("D3").Value = 1384.462953566667
$ws.Range("E3").Value = 4056.7215283

$ws.Range("B4").Value = 22887.87892265001
$ws.Range("D4").Value = 1476.7636146
$ws.Range("E4").Value = 4323.93000225

$ws.Range("B5").Value = 21931.65048931667
$ws.Range("D5").Value = 1458.52221285
$ws.Range("E5").Value = 4193.957928666667

$ws.Range("B6").Value = 22737.3828107
$ws.Range("D6").Value = 1467.364552866667
$ws.Range("E6").Value = 4417.14093795

$ws.Range("B7").Value = 21984.47007251667
$ws.Range("D7").Value = 1459.942688383333
$ws.Range("E7").Value = 4230.3473019

$ws.Range("B8").Value = 22742.92483245
$ws.Range("D8").Value = 1489.62160025
$ws.Range("E8").Value = 4276.627464116666

$ws.Range("B9").Value = 22881.75234935
$ws.Range("D9").Value = 1497.906630416667
$ws.Range("E9").Value = 4245.331106783334

$ws.Range("B10").Value = 22136.64730365001
$ws.Range("D10").Value = 1434.95356905
$ws.Range("E10").Value = 3920.31857385

$ws.Range("B11").Value = 22975.86661706667
$ws.Range("D11").Value = 1506.887080933333
$ws.Range("E11").Value = 4234.032952883334

$ws.Range("B12").Value = 22001.24391711667
$ws.Range("D12").Value = 1427.3509301
$ws.Range("E12").Value = 4024.42531745

$ws.Range("B13").Value = 22042.7830728
$ws.Range("D13").Value = 1431.153082666666
$ws.Range("E13").Value = 4176.7295289

